$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7180
$ws.Range("C3").Value = 160837
$ws.Range("C4").Value = 151861
$ws.Range("C8").Value = 64.40000000000001
